# Auto-generated edit script: update market-price derived columns (H-N)
# across the leve-profit tracking sheets, per scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 960.63635
$ws.Range("I18").Value = 638.5
$ws.Range("J18").Value = 1819.6666
$ws.Range("K18").Value = 638.5
$ws.Range("L18").Value = 1819.6666
$ws.Range("M18").Value = -354.5
$ws.Range("N18").Value = -2387.6666
$ws.Range("H33").Value = 620
$ws.Range("I33").Value = 345
$ws.Range("K33").Value = 345
$ws.Range("M33").Value = -116
$ws.Range("H40").Value = 4246.7
$ws.Range("I40").Value = 4053.4
$ws.Range("K40").Value = 4053.4
$ws.Range("M40").Value = -3878.4
$ws.Range("H103").Value = 493.29413
$ws.Range("J103").Value = 475.84616
$ws.Range("L103").Value = 1427.53848
$ws.Range("N103").Value = -2599.53848
$ws.Range("H137").Value = 3995.611
$ws.Range("I137").Value = 3175.889
$ws.Range("K137").Value = 9527.667000000001
$ws.Range("M137").Value = -6977.667000000001
$ws.Range("H138").Value = 2039.0857
$ws.Range("J138").Value = 1815.4286
$ws.Range("L138").Value = 5446.2858
$ws.Range("N138").Value = -15726.2858

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 39277.5
$ws.Range("I28").Value = 39277.5
$ws.Range("K28").Value = 39277.5
$ws.Range("M28").Value = -39085.5
$ws.Range("H32").Value = 6145.7554
$ws.Range("I32").Value = 2083.2974
$ws.Range("K32").Value = 2083.2974
$ws.Range("M32").Value = -1796.2974
$ws.Range("H41").Value = 18006.334
$ws.Range("I41").Value = 11019.667
$ws.Range("J41").Value = 24993
$ws.Range("K41").Value = 11019.667
$ws.Range("L41").Value = 24993
$ws.Range("M41").Value = -10605.667
$ws.Range("N41").Value = -25821
$ws.Range("H61").Value = 4512.615
$ws.Range("I61").Value = 3886.3333
$ws.Range("K61").Value = 3886.3333
$ws.Range("M61").Value = -3674.3333
$ws.Range("H74").Value = 1353.3334
$ws.Range("I74").Value = 778.5333000000001
$ws.Range("K74").Value = 778.5333000000001
$ws.Range("M74").Value = 95.46669999999995
$ws.Range("H77").Value = 1353.3334
$ws.Range("I77").Value = 778.5333000000001
$ws.Range("K77").Value = 3892.6665
$ws.Range("M77").Value = 475.3334999999997
$ws.Range("H97").Value = 877.8
$ws.Range("I97").Value = 877.8
$ws.Range("K97").Value = 877.8
$ws.Range("M97").Value = -381.8
$ws.Range("H99").Value = 39277.5
$ws.Range("I99").Value = 39277.5
$ws.Range("K99").Value = 39277.5
$ws.Range("M99").Value = -36282.5
$ws.Range("H110").Value = 2849.75
$ws.Range("I110").Value = 801.25
$ws.Range("K110").Value = 801.25
$ws.Range("M110").Value = 1243.75
$ws.Range("H132").Value = 2065.6128
$ws.Range("I132").Value = 1350.1
$ws.Range("K132").Value = 4050.3
$ws.Range("M132").Value = -1520.3
$ws.Range("H136").Value = 4512.615
$ws.Range("I136").Value = 3886.3333
$ws.Range("K136").Value = 11658.9999
$ws.Range("M136").Value = -9108.999899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5152.8335
$ws.Range("J31").Value = 4628
$ws.Range("L31").Value = 4628
$ws.Range("N31").Value = -5218
$ws.Range("H34").Value = 5152.8335
$ws.Range("J34").Value = 4628
$ws.Range("L34").Value = 4628
$ws.Range("N34").Value = -5032
$ws.Range("H94").Value = 6741.3335
$ws.Range("I94").Value = 6741.3335
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 6741.3335
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -6290.3335
$ws.Range("N94").ClearContents()
$ws.Range("H134").Value = 3116.1765
$ws.Range("I134").Value = 3409.5557
$ws.Range("K134").Value = 10228.6671
$ws.Range("M134").Value = -7693.667099999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 339.16666
$ws.Range("I8").Value = 339.16666
$ws.Range("K8").Value = 1017.49998
$ws.Range("M8").Value = -878.4999799999999
$ws.Range("H39").Value = 5414
$ws.Range("J39").Value = 8327
$ws.Range("L39").Value = 24981
$ws.Range("N39").Value = -25569
$ws.Range("H107").Value = 833.6111
$ws.Range("J107").Value = 796.875
$ws.Range("L107").Value = 2390.625
$ws.Range("N107").Value = -6230.625
$ws.Range("H121").Value = 91635
$ws.Range("I121").Value = 200078
$ws.Range("J121").Value = 1265.8334
$ws.Range("K121").Value = 600234
$ws.Range("L121").Value = 3797.5002
$ws.Range("M121").Value = -598924
$ws.Range("N121").Value = -6417.5002
$ws.Range("H129").Value = 4544
$ws.Range("I129").Value = 1108
$ws.Range("J129").Value = 7980
$ws.Range("K129").Value = 3324
$ws.Range("L129").Value = 23940
$ws.Range("M129").Value = 1676
$ws.Range("N129").Value = -33940

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3919.625
$ws.Range("I40").Value = 2753.7
$ws.Range("K40").Value = 2753.7
$ws.Range("M40").Value = -2617.7
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H122").Value = 6246.75
$ws.Range("I122").Value = 3990
$ws.Range("K122").Value = 11970
$ws.Range("M122").Value = -9520
$ws.Range("H132").Value = 2604.8235
$ws.Range("I132").Value = 2056.5454
$ws.Range("K132").Value = 6169.6362
$ws.Range("M132").Value = -3639.6362
$ws.Range("H136").Value = 2197.0908
$ws.Range("I136").Value = 1259.1666
$ws.Range("K136").Value = 3777.4998
$ws.Range("M136").Value = -1227.4998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 15999.667
$ws.Range("J10").Value = 3999.5
$ws.Range("L10").Value = 3999.5
$ws.Range("N10").Value = -4337.5
$ws.Range("H24").Value = 333348480
$ws.Range("J24").Value = 500002720
$ws.Range("L24").Value = 500002720
$ws.Range("N24").Value = -500003180
$ws.Range("H80").Value = 64999.332
$ws.Range("J80").Value = 64999.332
$ws.Range("L80").Value = 64999.332
$ws.Range("N80").Value = -66995.33199999999
$ws.Range("H83").Value = 64999.332
$ws.Range("J83").Value = 64999.332
$ws.Range("L83").Value = 194997.996
$ws.Range("N83").Value = -204981.996
$ws.Range("H122").Value = 3217.1177
$ws.Range("I122").Value = 2788.5
$ws.Range("K122").Value = 8365.5
$ws.Range("M122").Value = -5915.5
$ws.Range("H132").Value = 3515.1538
$ws.Range("I132").Value = 3780.5173
$ws.Range("J132").Value = 2745.6
$ws.Range("K132").Value = 11341.5519
$ws.Range("L132").Value = 8236.799999999999
$ws.Range("N132").Value = -13296.8
